$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store plain text values (e.g. "575.44",
# "  -1.18%  ") even though many look numeric. Force the whole D:E data range to
# Text format first so that assigning these strings does not get auto-converted
# into numeric cell values / floating point numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '64.026.79'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '2.756.49'
$ws.Range("E3").Value = '  -0.79%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '575.44'
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("D6").Value = '159.44'
$ws.Range("E6").Value = '  -1.22%  '
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").Value = '0.602'
$ws.Range("E8").Value = '  -3.42%  '
$ws.Range("E9").Value = '  -2.80%  '
$ws.Range("B10").Value = 'TRON'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D10").Value = '0.165'
$ws.Range("E10").Value = '  +4.06%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").Value = '5.78'
$ws.Range("E11").Value = '  -15.52%  '
$ws.Range("D12").Value = '0.388'
$ws.Range("E12").Value = '  -2.47%  '
$ws.Range("D13").Value = '3.244.94'
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("D14").Value = '27.02'
$ws.Range("E14").Value = '  -3.22%  '
$ws.Range("D15").Value = '63.657.93'
$ws.Range("E15").Value = '  -0.31%  '
$ws.Range("D16").Value = '0.0000151'
$ws.Range("E16").Value = '  -4.82%  '
$ws.Range("D17").Value = '2.761.93'
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").Value = '12.17'
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").Value = '4.86'
$ws.Range("E19").Value = '  -2.71%  '
$ws.Range("D20").Value = '356.66'
$ws.Range("E20").Value = '  -2.89%  '
$ws.Range("D21").Value = '6.70'
$ws.Range("E21").Value = '  -5.50%  '
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").Value = '0.529'
$ws.Range("E23").Value = '  -3.82%  '
$ws.Range("D24").Value = '65.24'
$ws.Range("E24").Value = '  -3.33%  '
$ws.Range("D25").Value = '0.170'
$ws.Range("E25").Value = '  -1.83%  '
$ws.Range("D26").Value = '8.58'
$ws.Range("E26").Value = '  -0.93%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("D28").Value = '0.0₃0905'
$ws.Range("E28").Value = '  -4.62%  '
$ws.Range("D29").Value = '7.30'
$ws.Range("E29").Value = '  +0.30%  '
$ws.Range("D30").Value = '1.94'
$ws.Range("E30").Value = '  -4.34%  '
$ws.Range("D31").Value = '1.25'
$ws.Range("E31").Value = '  -0.98%  '
$ws.Range("D32").Value = '169.13'
$ws.Range("E32").Value = '  -2.72%  '
$ws.Range("D33").Value = '4.94'
$ws.Range("E33").Value = '  -1.58%  '
$ws.Range("D34").Value = '20.15'
$ws.Range("E34").Value = '  -3.32%  '
$ws.Range("D35").Value = '1.48'
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").Value = '1.80'
$ws.Range("E37").Value = '  -2.43%  '
$ws.Range("E38").Value = '  -2.05%  '
$ws.Range("D39").Value = '351.92'
$ws.Range("E39").Value = '  +2.48%  '
$ws.Range("D40").Value = '6.27'
$ws.Range("E40").Value = '  +0.38%  '
$ws.Range("D41").Value = '4.18'
$ws.Range("E41").Value = '  -2.34%  '
$ws.Range("D42").Value = '39.06'
$ws.Range("E42").Value = '  -1.79%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = '21.93'
$ws.Range("E43").Value = '  -3.06%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '21.43'
$ws.Range("E44").Value = '  -2.94%  '
$ws.Range("D45").Value = '0.0589'
$ws.Range("E45").Value = '  -3.05%  '
$ws.Range("D46").Value = '0.633'
$ws.Range("E46").Value = '  -2.83%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '136.12'
$ws.Range("E47").Value = '  -1.41%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '0.0254'
$ws.Range("E48").Value = '  -2.63%  '
$ws.Range("E49").Value = '  -1.00%  '
$ws.Range("E50").Value = '  +0.38%  '
$ws.Range("D51").Value = '11.04'
$ws.Range("E51").Value = '  +0.09%  '

# Restore the original (default/no explicit number format) cell formatting by
# pasting the format from a cell that was never touched and always kept the
# workbook's default style, so the resulting cells look exactly like the rest
# of the untouched text cells (no stray "Text" number format left behind).
$ws.Range("B2").Copy() | Out-Null
$dataRange.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
